# Support for any other fields
# Adds a new "Priority" column to the test-case table, populates it for the
# three test-case header rows, fills in a couple of other fields that were
# left blank (Area Path / State) and cleans up unused styling left over on
# empty helper cells in the "Description"/"Automation Status"... area
# (columns L/M) of the step rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Add the new "Priority" table column (grows the table from M to N)
# ---------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.ListColumns.Add() | Out-Null
$ws.Cells.Item(1, 14).Value = "Priority"

# ---------------------------------------------------------------------
# 2. Fill in the Priority values for the three test case rows
# ---------------------------------------------------------------------
$ws.Range("N2").Value = 2
$ws.Range("N6").Value = 3
$ws.Range("N10").Value = 1

# ---------------------------------------------------------------------
# 3. Fill in a couple of other fields that were previously empty
# ---------------------------------------------------------------------
$ws.Range("G2").Value = "specsync-plugins-demo\Area2\Area 2.1"
$ws.Range("J6").Value = "tag1"
$ws.Range("G10").Value = "specsync-plugins-demo"
$ws.Range("I10").Value = "Ready"

# ---------------------------------------------------------------------
# 4. Clean up the left-over styling in columns L/M:
#    - cells that already hold a value just lose the stray style
#    - cells that are blank are cleared out completely
# ---------------------------------------------------------------------
$valuedCells = @("L2", "L6", "M6", "L10", "M10")
foreach ($c in $valuedCells) {
    $ws.Range($c).ClearFormats()
}

$blankCells = @(
    "M2",
    "L3", "M3",
    "L4", "M4",
    "L5", "M5",
    "L7", "M7",
    "L8", "M8",
    "L9", "M9",
    "L11", "M11",
    "L12", "M12"
)
foreach ($c in $blankCells) {
    $ws.Range($c).Clear()
}

# ---------------------------------------------------------------------
# 5. Restore the selection to match the active cell used while editing
# ---------------------------------------------------------------------
$ws.Range("I10").Select() | Out-Null
